# Updated cryptos list on Mon May 29 07:27:17 UTC 2023 with GitHub Actions
#
# Applies the refreshed price/volume snapshot to the cryptos table on the
# active sheet. The "Price" column (D) holds values that look numeric
# (e.g. "1.006", "28.061.68") but are authored as literal text, so a plain
# ".Value = ..." assignment would let Excel's COM layer silently reinterpret
# them as numbers. We force literal text by priming NumberFormat to "@"
# before the write, then clear the now-unneeded explicit format afterwards
# so the cell's style/format stays byte-for-byte what it was before (no
# lingering "@" number format attached to the cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($sheet, $addr, $val) {
    $r = $sheet.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

# Row 2 - Bitcoin
Set-TextValue $ws "D2" "28.061.68"
$ws.Range("E2").Value = "  +2.03%  "

# Row 3 - Ethereum
Set-TextValue $ws "D3" "1.911.36"
$ws.Range("E3").Value = "  +2.31%  "

# Row 4 - TetherUSD
Set-TextValue $ws "D4" "1.006"
$ws.Range("E4").Value = "  -0.79%  "

# Row 5 - BNB
Set-TextValue $ws "D5" "317.00"
$ws.Range("E5").Value = "  +1.47%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  -0.81%  "

# Row 7 - XRP
Set-TextValue $ws "D7" "0.4822"
$ws.Range("E7").Value = "  +0.88%  "

# Row 8 - Cardano
Set-TextValue $ws "D8" "0.3816"
$ws.Range("E8").Value = "  +1.02%  "

# Row 9 - Dogecoin
Set-TextValue $ws "D9" "0.07365"
$ws.Range("E9").Value = "  +0.07%  "

# Row 10 - Polygon
Set-TextValue $ws "D10" "0.9346"
$ws.Range("E10").Value = "  -0.22%  "

# Row 11 - Solana
$ws.Range("E11").Value = "  +0.49%  "

# Row 12 - TRON
Set-TextValue $ws "D12" "0.07787"
$ws.Range("E12").Value = "  -0.96%  "

# Row 13 - WrappedEther
Set-TextValue $ws "D13" "1.936.40"
$ws.Range("E13").Value = "  +3.61%  "

# Row 14 - Polkadot
Set-TextValue $ws "D14" "5.513"
$ws.Range("E14").Value = "  +1.38%  "

# Row 15 - Chainlink
Set-TextValue $ws "D15" "6.633"
$ws.Range("E15").Value = "  +0.91%  "

# Row 16 - Litecoin
Set-TextValue $ws "D16" "91.83"
$ws.Range("E16").Value = "  +1.23%  "

# Row 17 - BinanceUSD
$ws.Range("E17").Value = "  -0.77%  "

# Row 18 - ShibaInu
Set-TextValue $ws "D18" "0.000008840"
$ws.Range("E18").Value = "  -0.67%  "

# Row 19 - Dai
$ws.Range("E19").Value = "  -0.76%  "

# Row 20 - WrappedBTC
Set-TextValue $ws "D20" "28.085.36"
$ws.Range("E20").Value = "  +2.05%  "

# Row 21 - Avalanche
Set-TextValue $ws "D21" "14.84"
$ws.Range("E21").Value = "  -0.52%  "

# Row 22 - Uniswap
Set-TextValue $ws "D22" "5.182"
$ws.Range("E22").Value = "  +0.98%  "

# Row 23 - WrappedliquidstakedEther2.0
Set-TextValue $ws "D23" "2.135.64"
$ws.Range("E23").Value = "  +1.41%  "

# Row 24 - Cosmos
$ws.Range("E24").Value = "  +1.93%  "

# Row 25 - ...
Set-TextValue $ws "D25" "155.69"
$ws.Range("E25").Value = "  +1.20%  "

# Row 26
$ws.Range("E26").Value = "  -1.89%  "

# Row 27
$ws.Range("E27").Value = "  +0.22%  "

# Row 28 - LidoDAOToken
Set-TextValue $ws "D28" "2.122"
$ws.Range("E28").Value = "  +4.90%  "

# Row 29 - BitcoinCash
Set-TextValue $ws "D29" "116.72"
$ws.Range("E29").Value = "  +0.67%  "

# Row 30 - InternetComputer(DFINITY)
Set-TextValue $ws "D30" "4.963"
$ws.Range("E30").Value = "  -0.75%  "

# Row 31 - Stellar
Set-TextValue $ws "D31" "0.08963"
$ws.Range("E31").Value = "  +0.39%  "

# Row 32 - HuobiToken
Set-TextValue $ws "D32" "3.312"
$ws.Range("E32").Value = "  -0.74%  "

# Row 33 - ARBITRUM
Set-TextValue $ws "D33" "1.255"
$ws.Range("E33").Value = "  +3.20%  "

# Row 34 - ImmutableX
Set-TextValue $ws "D34" "0.7785"
$ws.Range("E34").Value = "  +3.37%  "

# Row 35 - Filecoin
$ws.Range("E35").Value = "  +1.44%  "

# Row 36 - RenderToken
Set-TextValue $ws "D36" "2.662"
$ws.Range("E36").Value = "  -1.33%  "

# Row 37 - VeChain
Set-TextValue $ws "D37" "0.02056"
$ws.Range("E37").Value = "  +0.21%  "

# Row 38 - TrustWalletToken
Set-TextValue $ws "D38" "1.111"
$ws.Range("E38").Value = "  -0.80%  "

# Row 39 - Hedera
Set-TextValue $ws "D39" "0.05322"
$ws.Range("E39").Value = "  +0.87%  "

# Row 40 - TheSandbox
Set-TextValue $ws "D40" "0.5482"
$ws.Range("E40").Value = "  +2.60%  "

# Row 41 - MXToken
Set-TextValue $ws "D41" "2.988"
$ws.Range("E41").Value = "  -0.57%  "

# Row 42 - FraxShare
Set-TextValue $ws "D42" "7.024"
$ws.Range("E42").Value = "  -0.86%  "

# Row 43 & 44 - Algorand/Aptos swapped rank order
Set-TextValue $ws "B43" "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws "D43" "8.539"
$ws.Range("E43").Value = "  +0.74%  "

Set-TextValue $ws "B44" "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws "D44" "0.1529"
$ws.Range("E44").Value = "  +0.23%  "

# Row 45 - EnergySwap
Set-TextValue $ws "D45" "10.74"
$ws.Range("E45").Value = "  +0.67%  "

# Row 46 - Decentraland
Set-TextValue $ws "D46" "0.4843"
$ws.Range("E46").Value = "  +0.84%  "

# Row 47 - Quant
Set-TextValue $ws "D47" "108.56"
$ws.Range("E47").Value = "  +5.59%  "

# Row 48 - PaxDollar
$ws.Range("E48").Value = "  -0.82%  "

# Row 49 - NEARProtocol
Set-TextValue $ws "D49" "1.656"
$ws.Range("E49").Value = "  -0.17%  "

# Row 50 - Aave
Set-TextValue $ws "D50" "68.12"
$ws.Range("E50").Value = "  +0.94%  "

# Row 51 - Cronos
Set-TextValue $ws "D51" "0.06087"
$ws.Range("E51").Value = "  +0.00%  "
